$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right 6 -> 9, Wrong 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 (Total): Right 66 -> 99, Wrong -6 -> -4, Max fraction 60/168 -> 95/252
$ws.Range("B12").Value = 99
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "95/252"
